# Applies the "Getters / Setters" rows to the Funcoes workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the order the author typed them: columns A and B for both new
# rows first (Getters / Setters pair), then column C for both rows - this
# matches the shared-string insertion order recorded in the workbook.
$ws.Range("A34").Value = "Getters"
$ws.Range("B34").Value = "Onde eu vou buscar valor"
$ws.Range("A35").Value = "Setters "
$ws.Range("B35").Value = "Onde eu vou incluir valor"
$ws.Range("C34").Value = "Onde eu vou ler o valor da variavel"
$ws.Range("C35").Value = "Onde eu posso modificar a variavel"

# Apply the same formatting used by the rest of the data rows (style index 1:
# size-18 font, centered, wrap text), matching row 33's formatting, by
# copying the row's format rather than building a brand new style.
$ws.Range("A33:C33").Copy() | Out-Null
$ws.Range("A34:C35").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the view to scroll down to the newly added rows, similar to the
# author's resulting selection/scroll position.
$ws.Range("C36").Select()
$excel.ActiveWindow.ScrollRow = 26
